# Laporan Buku Tamu - menambahkan fitur upload gambar
# (adds new guest-book rows + fixes a typo in an existing row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing typo: "perong" -> "peron" (row 2, NAMA TAMU) ---
$ws.Range("C2").Value = "peron"

# Helper: write a value that must stay plain TEXT (no Excel auto type
# inference into a date/number) without leaving behind any cell-style
# residue. Going through Formula -> Copy -> PasteSpecial(xlPasteValues)
# keeps the result a literal shared string with no number-format churn.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# --- row 7 ---
$ws.Range("A7").Value = 6
Set-TextValue $ws.Range("B7") "2024-10-14"
$ws.Range("C7").Value = "ival tim"
$ws.Range("D7").Value = "leles"
Set-TextValue $ws.Range("E7") "08578273282"
$ws.Range("F7").Value = "egi tim"
$ws.Range("G7").Value = "ngomongin tim"

# --- row 8 ---
$ws.Range("A8").Value = 7
Set-TextValue $ws.Range("B8") "2024-10-14"
$ws.Range("C8").Value = "test"
$ws.Range("D8").Value = "tekjtasldk"
$ws.Range("E8").Value = 1232131
$ws.Range("F8").Value = "skdfsdjf"
$ws.Range("G8").Value = "dfjsadfhk"

# --- row 9 ---
$ws.Range("A9").Value = 8
Set-TextValue $ws.Range("B9") "2024-10-15"
$ws.Range("C9").Value = "test2"
$ws.Range("D9").Value = "sdfkjasdjf"
$ws.Range("E9").Value = "ksdfakdsj"
$ws.Range("F9").Value = "kdfkaj"
$ws.Range("G9").Value = "kknsdflak"

# --- row 10 ---
$ws.Range("A10").Value = 9
Set-TextValue $ws.Range("B10") "2024-10-17"
$ws.Range("C10").Value = "azid"
$ws.Range("D10").Value = "blk"
$ws.Range("E10").Value = 1242432
$ws.Range("F10").Value = "turky"
$ws.Range("G10").Value = "ngendong"

# --- row 11 ---
$ws.Range("A11").Value = 10
Set-TextValue $ws.Range("B11") "2024-10-22"
$ws.Range("C11").Value = "bebey"
$ws.Range("D11").Value = "gg guntur"
$ws.Range("E11").Value = 66587970854
$ws.Range("F11").Value = "ffsdfg"
$ws.Range("G11").Value = "unity"
